$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.180.00'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '2.305.84'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''310.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''101.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.01%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '  +2.27%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.524'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.42%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''35.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.62%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.0811'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.94%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.112'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''6.97'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.93%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '2.662.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.12%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''14.95'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '2.303.77'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.06%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''0.809'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.68%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '43.092.59'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.85%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''12.48'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '0.0₃0918'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''6.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.09%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''68.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''240.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''2.04'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.84%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''2.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.85%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('B27').Value = 'LEO'
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').Value = '''3.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').Value = '''24.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.62%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = '''37.44'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.04%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = '''9.60'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = '''2.11'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = 'Monero'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = '''167.94'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.75%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = '''5.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = '''1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''17.78'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''0.0743'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.24%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = '''3.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''2.43'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.69%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = '''0.107'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.60%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = '''1.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.13%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = '''4.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.91%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = '''19.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.17%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = '''2.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.31%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '1.972.54'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = 'VeChain'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''0.0289'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.02%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''2.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.28%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''9.76'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = '''3.01'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +19.85%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = '''56.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.94%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = '2.531.85'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.98%  '
$ws.Range('E51').Style = 'Normal'
